# More robust metadata refresh url extraction.
# Inserts a new row (Moormerland) as row 2 of the BLP-URLs sheet, shifting
# all following rows down by one, and updates the _FilterDatabase defined
# name range to match the new used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 ("Bad Iburg,Stadt"),
# pushing all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# Populate the new row with the Moormerland entry.
$ws.Range("A2").Value = 457014
$ws.Range("B2").Value = "Moormerland"
$ws.Range("C2").Value = 53.314314
$ws.Range("D2").Value = 7.485564
$ws.Range("E2").Value = "http://www.mmld.de/download"
$ws.Range("F2").Value = "http://lkleer.maps.arcgis.com/home/webmap/viewer.html?webmap=e4311f176259429d970921af4cf49ab2"

# Match the original sheet's row-height convention for this row (13.8).
$ws.Rows.Item(2).RowHeight = 13.8

# The sheet's _FilterDatabase defined name pinned the data range to
# $A$1:$G$406; after the inserted row it must cover one more row.
$name = $wb.Names.Item("_xlnm._FilterDatabase")
$name.RefersTo = "='BLP-URLs'!`$A`$1:`$G`$407"

Write-Output "Inserted Moormerland row; filter database range extended to G407"
